# Update the "取得日時" (acquired-at) timestamp for all data rows on the
# "ランサーズ" sheet from 2026-02-10 07:01:11 to 2026-02-10 07:10:42.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

$oldTimestamp = "2026-02-10 07:01:11"
$newTimestamp = "2026-02-10 07:10:42"

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 1)
    if ($cell.Value2 -eq $oldTimestamp) {
        $cell.Value = $newTimestamp
    }
}
